# Updates the cryptos price/volume snapshot (and fixes the EnergySwap /
# Decentraland row ordering) to match the latest scrape.
#
# Note: several "Price" values look like plain numbers (e.g. 1.004, 0.4600).
# Assigning them straight to .Value would make Excel auto-convert them to
# real numbers (losing trailing zeros / exact decimal text). To keep them
# as text - matching the original inline-string cells - those values are
# prefixed with a leading apostrophe, Excel's standard "treat as text"
# quote-prefix marker, which is stripped from the stored value.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range('D2').Value = '30.472.97'
$ws.Range('E2').Value = '  -0.88%  '
$ws.Range('D3').Value = '2.108.30'
$ws.Range('E3').Value = '  -0.03%  '
$ws.Range('D4').Value = '''1.004'
$ws.Range('E4').Value = '  +0.29%  '
$ws.Range('D5').Value = '''333.85'
$ws.Range('E5').Value = '  +0.20%  '
$ws.Range('D6').Value = '''1.003'
$ws.Range('E6').Value = '  +0.32%  '
$ws.Range('D7').Value = '''0.5245'
$ws.Range('E7').Value = '  -0.90%  '
$ws.Range('D8').Value = '''0.4600'
$ws.Range('E8').Value = '  +5.66%  '
$ws.Range('D9').Value = '''53.68'
$ws.Range('E9').Value = '  +13.56%  '
$ws.Range('D10').Value = '''0.08966'
$ws.Range('E10').Value = '  +0.07%  '
$ws.Range('E11').Value = '  +0.88%  '
$ws.Range('E12').Value = '  -1.30%  '
$ws.Range('D13').Value = '2.104.88'
$ws.Range('E13').Value = '  -0.02%  '
$ws.Range('D14').Value = '''6.784'
$ws.Range('E14').Value = '  +0.88%  '
$ws.Range('D15').Value = '''7.839'
$ws.Range('E15').Value = '  +0.97%  '
$ws.Range('D16').Value = '''96.56'
$ws.Range('E16').Value = '  -0.23%  '
$ws.Range('E17').Value = '  +0.36%  '
$ws.Range('D18').Value = '''0.00001131'
$ws.Range('E18').Value = '  -0.04%  '
$ws.Range('D19').Value = '''0.06629'
$ws.Range('E19').Value = '  -0.86%  '
$ws.Range('D20').Value = '''19.26'
$ws.Range('D21').Value = '''1.003'
$ws.Range('E21').Value = '  +0.15%  '
$ws.Range('D22').Value = '''6.284'
$ws.Range('E22').Value = '  -0.30%  '
$ws.Range('D23').Value = '30.543.42'
$ws.Range('E23').Value = '  -0.86%  '
$ws.Range('E24').Value = '  +0.71%  '
$ws.Range('D25').Value = '''2.359'
$ws.Range('E25').Value = '  +3.40%  '
$ws.Range('D26').Value = '2.353.69'
$ws.Range('E26').Value = '  +0.08%  '
$ws.Range('D27').Value = '''22.32'
$ws.Range('E27').Value = '  -1.14%  '
$ws.Range('D28').Value = '''2.564'
$ws.Range('E28').Value = '  -0.48%  '
$ws.Range('D29').Value = '''163.49'
$ws.Range('E29').Value = '  +0.77%  '
$ws.Range('D30').Value = '''132.77'
$ws.Range('E30').Value = '  -0.09%  '
$ws.Range('D31').Value = '''1.196'
$ws.Range('E31').Value = '  +0.35%  '
$ws.Range('D32').Value = '''0.1072'
$ws.Range('E32').Value = '  -0.78%  '
$ws.Range('D33').Value = '''1.690'
$ws.Range('E33').Value = '  +9.17%  '
$ws.Range('D34').Value = '''6.152'
$ws.Range('E34').Value = '  -0.28%  '
$ws.Range('D35').Value = '''3.934'
$ws.Range('E35').Value = '  +0.99%  '
$ws.Range('D36').Value = '''10.45'
$ws.Range('E36').Value = '  +9.00%  '
$ws.Range('D37').Value = '''0.02573'
$ws.Range('E37').Value = '  -0.72%  '
$ws.Range('D38').Value = '''0.06821'
$ws.Range('E38').Value = '  +0.93%  '
$ws.Range('D39').Value = '''5.548'
$ws.Range('E39').Value = '  +0.23%  '
$ws.Range('D40').Value = '''12.81'
$ws.Range('E40').Value = '  +1.43%  '
$ws.Range('D41').Value = '''0.2289'
$ws.Range('E41').Value = '  +0.80%  '
$ws.Range('D42').Value = '''0.6886'
$ws.Range('E42').Value = '  +0.75%  '
$ws.Range('E43').Value = '  +0.05%  '
$ws.Range('D44').Value = '''2.350'
$ws.Range('E44').Value = '  +5.73%  '
$ws.Range('E45').Value = '  +0.29%  '
$ws.Range('B46').Value = 'Decentraland'
$ws.Range('C46').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D46').Value = '''0.6374'
$ws.Range('E46').Value = '  -0.62%  '
$ws.Range('B47').Value = 'EnergySwap'
$ws.Range('C47').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D47').Value = '''13.93'
$ws.Range('E47').Value = '  -0.73%  '
$ws.Range('E48').Value = '  +0.00%  '
$ws.Range('D49').Value = '''0.00000000349'
$ws.Range('E49').Value = '  +22.92%  '
$ws.Range('D50').Value = '''1.245'
$ws.Range('E50').Value = '  -1.17%  '
$ws.Range('D51').Value = '''1.222'
$ws.Range('E51').Value = '  +2.70%  '
